{"js": "// The document has two consecutive empty paragraphs right after\n// \"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\" (and right before the \"16.06.2023\" paragraph).\n// The edit collapses those two empty paragraphs into a single empty\n// paragraph that carries Word's \"_GoBack\" bookmark (the marker Word drops\n// at the location of the last edit).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two consecutive empty paragraphs that sit between\n// \"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\" and \"16.06.2023\".\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  if (\n    paragraphs.items[i].text === \"\" &&\n    paragraphs.items[i + 1].text === \"\" &&\n    i > 0 &&\n    paragraphs.items[i - 1].text.indexOf(\"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\") !== -1\n  ) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the pair of empty paragraphs to merge.\");\n}\n\nconst firstEmpty = paragraphs.items[targetIndex];\nconst secondEmpty = paragraphs.items[targetIndex + 1];\n\n// Remove the second (now redundant) empty paragraph...\nsecondEmpty.delete();\n\n// ...and drop the \"_GoBack\" bookmark into the remaining empty paragraph,\n// matching what Word stamps at the last edited / cursor location.\nfirstEmpty.getRange().insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The document has two consecutive empty paragraphs right after the\n# \"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\" paragraph (and right before the \"16.06.2023\"\n# paragraph). The edit collapses those two empty paragraphs into a single\n# empty paragraph that carries Word's \"_GoBack\" bookmark - the marker Word\n# drops at the location of the last edit / cursor position.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\") so the edit is\n# not dependent on a hard-coded paragraph index.\n$anchor = $d.Content\n$found = $anchor.Find.Execute(\"G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131\")\nif (-not $found -or -not $anchor.Find.Found) {\n    throw \"Could not find the 'G\u00f6rev da\u011f\u0131l\u0131m\u0131 yap\u0131ld\u0131' paragraph.\"\n}\n\n$anchorIndex = $anchor.Paragraphs.Item(1).Index\n$firstEmpty = $d.Paragraphs.Item($anchorIndex + 1)\n$secondEmpty = $d.Paragraphs.Item($anchorIndex + 2)\n\nif ($firstEmpty.Range.Text.Trim() -ne \"\" -or $secondEmpty.Range.Text.Trim() -ne \"\") {\n    throw \"Expected two empty paragraphs after the anchor paragraph.\"\n}\n\n# Remove the second (now redundant) empty paragraph...\n$secondEmpty.Range.Delete()\n\n# ...and drop the \"_GoBack\" bookmark into the remaining empty paragraph.\n$d.Bookmarks.Add(\"_GoBack\", $firstEmpty.Range)\n"}
